# Regenerate save_data: column G ("K" = strikeouts) values recomputed
# and rewritten for rows 2-29 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 7
    3  = 5
    4  = 5
    5  = 9
    6  = 4
    7  = 4
    8  = 5
    9  = 8
    10 = 6
    11 = 3
    12 = 6
    13 = 2
    14 = 5
    15 = 7
    16 = 4
    17 = 10
    18 = 2
    19 = 7
    20 = 7
    21 = 8
    22 = 7
    23 = 10
    24 = 5
    25 = 4
    26 = 3
    27 = 1
    28 = 2
    29 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
